# "Chrome and firefox issue sorted"
#
# The real, functional change behind this commit is the introduction of a
# new worksheet, "LoginPageUI", that holds the localisation keys / UI text
# used to verify the login page across browsers. It is inserted between
# the existing "InvalidUserNameLogin" and "InvalidPasswordLogin" sheets:
#
#   loginData, InvalidUserNameLogin, LoginPageUI, InvalidPasswordLogin
#
# New sheet layout:
#   A1: loginpagetitle            B1: headerlogonmetext            C1: verifytermsandconditions
#   A2: Log On - Ci Anywhere      B2: Log on using your details    C2: View Terms and Conditions
#
# (A2 re-uses the "Log On - Ci Anywhere" string already shared with the
# "loginData"/"InvalidUserNameLogin" sheets.)

$wb = $excel.ActiveWorkbook

# Leave the cursor on the first sheet parked at A2 (matches the saved
# selection left behind once the new sheet became the active tab).
$ws1 = $wb.Worksheets.Item("loginData")
$ws1.Range("A2").Select() | Out-Null

# Insert the new worksheet right before "InvalidPasswordLogin" so the tab
# order becomes loginData, InvalidUserNameLogin, LoginPageUI, InvalidPasswordLogin.
$wsBefore = $wb.Worksheets.Item("InvalidPasswordLogin")
$ws = $wb.Worksheets.Add($wsBefore)
$ws.Name = "LoginPageUI"

# Populate the new sheet. Values are written in this order so that any
# newly-introduced shared strings land in the same order they were first
# used (row 2 values, then row 1 headers, then the reused A2 value).
$ws.Range("B2").Value = "Log on using your details"
$ws.Range("C2").Value = "View Terms and Conditions"
$ws.Range("A1").Value = "loginpagetitle"
$ws.Range("B1").Value = "headerlogonmetext"
$ws.Range("C1").Value = "verifytermsandconditions"
$ws.Range("A2").Value = "Log On - Ci Anywhere"

# Leave the new sheet as the active tab/selection, parked where the author
# left it.
$ws.Range("E22").Select() | Out-Null
